$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''63.520.37'
$ws.Range("E2").Value = '  +2.53%  '
$ws.Range("D3").Value = '''3.127.88'
$ws.Range("E3").Value = '  +1.32%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '''587.81'
$ws.Range("E5").Value = '  +1.18%  '
$ws.Range("D6").Value = '''146.39'
$ws.Range("E6").Value = '  +2.76%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '''3.121.13'
$ws.Range("E8").Value = '  +1.45%  '
$ws.Range("E9").Value = '  +0.97%  '
$ws.Range("D10").Value = '''0.160'
$ws.Range("E10").Value = '  +13.71%  '
$ws.Range("E11").Value = '  +0.36%  '
$ws.Range("E12").Value = '  +0.33%  '
$ws.Range("E13").Value = '  +4.52%  '
$ws.Range("D14").Value = '''36.63'
$ws.Range("E14").Value = '  +3.49%  '
$ws.Range("E15").Value = '  -0.73%  '
$ws.Range("D16").Value = '''3.643.05'
$ws.Range("E16").Value = '  +1.24%  '
$ws.Range("E17").Value = '  -1.59%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '''63.452.36'
$ws.Range("E18").Value = '  +2.55%  '
$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = '''3.124.97'
$ws.Range("E19").Value = '  +1.35%  '
$ws.Range("D20").Value = '''463.54'
$ws.Range("E20").Value = '  +3.43%  '
$ws.Range("D21").Value = '''14.38'
$ws.Range("E21").Value = '  +3.05%  '
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("E23").Value = '  +1.31%  '
$ws.Range("E24").Value = '  -4.10%  '
$ws.Range("D25").Value = '''82.11'
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("E26").Value = '  -0.17%  '
$ws.Range("D27").Value = '''8.91'
$ws.Range("E27").Value = '  +8.25%  '
$ws.Range("E28").Value = '  +1.02%  '
$ws.Range("D29").Value = '''2.23'
$ws.Range("E29").Value = '  -1.73%  '
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("D31").Value = '''6.88'
$ws.Range("E31").Value = '  +1.59%  '
$ws.Range("D32").Value = '''27.07'
$ws.Range("E32").Value = '  +0.63%  '
$ws.Range("E33").Value = '  -2.38%  '
$ws.Range("D34").Value = '''0.0₃0874'
$ws.Range("E34").Value = '  +9.71%  '
$ws.Range("B35").Value = 'Mantle'
$ws.Range("C35").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D35").Value = '''1.05'
$ws.Range("E35").Value = '  +1.45%  '
$ws.Range("B36").Value = 'Stacks'
$ws.Range("C36").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D36").Value = '''2.35'
$ws.Range("E36").Value = '  +7.14%  '
$ws.Range("E37").Value = '  +13.32%  '
$ws.Range("D38").Value = '''6.10'
$ws.Range("E38").Value = '  +0.65%  '
$ws.Range("D39").Value = '''51.03'
$ws.Range("E39").Value = '  +1.33%  '
$ws.Range("D40").Value = '''449.94'
$ws.Range("E40").Value = '  +6.51%  '
$ws.Range("D41").Value = '''8.76'
$ws.Range("E41").Value = '  -0.63%  '
$ws.Range("D42").Value = '''0.0372'
$ws.Range("E42").Value = '  -0.15%  '
$ws.Range("D43").Value = '''2.899.26'
$ws.Range("E43").Value = '  +0.79%  '
$ws.Range("D44").Value = '''0.281'
$ws.Range("E44").Value = '  +2.69%  '
$ws.Range("E45").Value = '  +1.44%  '
$ws.Range("D46").Value = '''2.18'
$ws.Range("E46").Value = '  +1.84%  '
$ws.Range("D47").Value = '''36.38'
$ws.Range("E47").Value = '  +3.34%  '
$ws.Range("D48").Value = '''126.11'
$ws.Range("E48").Value = '  +1.52%  '
$ws.Range("E49").Value = '  +0.04%  '
$ws.Range("E50").Value = '  -0.20%  '
$ws.Range("D51").Value = '''24.79'
$ws.Range("E51").Value = '  +2.65%  '

# Reset style on price cells so the quote-prefix formatting introduced by the
# leading apostrophe does not leave a stray style index behind.
$priceRefs = @("D2","D3","D5","D6","D8","D10","D14","D16","D18","D19","D20","D21","D25","D27","D29","D31","D32","D34","D35","D36","D38","D39","D40","D41","D42","D43","D44","D46","D47","D48","D51")
foreach ($ref in $priceRefs) {
    $ws.Range($ref).Style = "Normal"
}
